$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.489.47'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.808.98'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '225.82'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = '0.599'
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '36.31'
$ws.Range("E8").Value = '  +3.71%  '
$ws.Range("D9").Value = '0.293'
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("D10").Value = '0.0681'
$ws.Range("E10").Value = '  -1.63%  '
$ws.Range("D11").Value = '0.0966'
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").Value = '2.069.42'
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").Value = '11.32'
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("D14").Value = '1.823.49'
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("D15").Value = '0.630'
$ws.Range("E15").Value = '  -1.68%  '
$ws.Range("D16").Value = '34.463.54'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").Value = '4.42'
$ws.Range("E17").Value = '  +1.84%  '
$ws.Range("D18").Value = '68.63'
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").Value = '242.95'
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").Value = '0.0₃0774'
$ws.Range("E20").Value = '  -2.86%  '
$ws.Range("D21").Value = '11.25'
$ws.Range("E21").Value = '  -2.23%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = '4.11'
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("D24").Value = '2.22'
$ws.Range("E24").Value = '  +5.14%  '
$ws.Range("D25").Value = '171.12'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").Value = '7.94'
$ws.Range("E26").Value = '  +3.52%  '
$ws.Range("D27").Value = '17.26'
$ws.Range("E27").Value = '  +3.14%  '
$ws.Range("D28").Value = '0.121'
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").Value = '3.82'
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").Value = '3.93'
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("D32").Value = '1.23'
$ws.Range("E32").Value = '  -1.07%  '
$ws.Range("D33").Value = '0.0517'
$ws.Range("E33").Value = '  -2.02%  '
$ws.Range("D34").Value = '1.80'
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("D35").Value = '1.362.03'
$ws.Range("E35").Value = '  -2.36%  '
$ws.Range("D36").Value = '0.653'
$ws.Range("E36").Value = '  -3.64%  '
$ws.Range("D37").Value = '1.06'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  -5.32%  '
$ws.Range("D39").Value = '0.0186'
$ws.Range("E39").Value = '  -1.79%  '
$ws.Range("D40").Value = '2.42'
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("D41").Value = '2.78'
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("D42").Value = '81.06'
$ws.Range("E42").Value = '  -2.46%  '
$ws.Range("D43").Value = '0.937'
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").Value = '1.17'
$ws.Range("E44").Value = '  +4.89%  '
$ws.Range("D45").Value = '13.41'
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("D46").Value = '0.0500'
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("D47").Value = '1.971.61'
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").Value = '5.84'
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").Value = '102.74'
$ws.Range("E50").Value = '  -1.81%  '
$ws.Range("D51").Value = '0.0₆0122'
$ws.Range("E51").Value = '  -5.87%  '
